$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the tech/description table for rows 5-22 (rows 1-4 are unchanged).
# A new row ("NG_CC_Existing") replaces the old row 5/6 pair, everything
# below ripples up by one, and the old DFO_CT_Existing/DFO_CA_Existing pair
# collapses into a single DFO_CC_Existing row — net effect: the table shrinks
# from 24 data rows down to 22, so the two now-unused trailing rows are
# deleted afterwards.

$data = @(
    @("NG_CC_Existing", "Natural Gas Fired Combined Cycle"),
    @("BLQ_ST_Existing", "Wood/Wood Waste Biomass"),
    @("SUN_PV_Existing", "Solar Photovoltaic"),
    @("MWH_BA1h_Existing", "Batteries"),
    @("DFO_GT_Existing", "Petroleum Liquids"),
    @("WDS_ST_Existing", "Wood/Wood Waste Biomass"),
    @("WH_ST_Existing", "All Other"),
    @("LFG_IC_Existing", "Landfill Gas"),
    @("WND_WT_Existing", "Onshore Wind Turbine"),
    @("AB_ST_Existing", "Other Waste Biomass"),
    @("NG_ST_Existing", "Natural Gas Steam Turbine"),
    @("WAT_HY_Existing", "Conventional Hydroelectric"),
    @("WAT_PS_Existing", "Hydroelectric Pumped Storage"),
    @("DFO_CC_Existing", "Petroleum Liquids"),
    @("BIT_ST_Existing", "Conventional Steam Coal"),
    @("LFG_GT_Existing", "Landfill Gas"),
    @("OBG_IC_Existing", "Other Waste Biomass"),
    @("MWH_BA2h_Existing", "Batteries")
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# The table now ends at row 22 (was row 24) - remove the two leftover rows.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(23).Delete()
